# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" worksheets, reflecting a regenerated gh-pages
# data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for sheet "展览"
$exhibitionUpdates = @{
    2  = 254
    12 = 700
    13 = 763
    14 = 1519
    15 = 1519
    20 = 330
    23 = 104
    24 = 6638
    25 = 5002
    29 = 204
    32 = 1287
    35 = 618
    38 = 251
}

# Row -> new F-column value for sheet "全部类型"
$allTypesUpdates = @{
    2  = 254
    16 = 700
    17 = 763
    18 = 1519
    19 = 1519
    24 = 330
    26 = 104
    29 = 6638
    30 = 5002
    32 = 204
    34 = 1287
    39 = 618
    43 = 251
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
